$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 68 ---------------------------------------------------------------
# Copy formatting from the last existing data row (67) onto the two new
# rows before writing values, so the new rows carry the same styles
# (bold/bordered index column, date-time number format column) as the
# rest of the table.
$ws.Range("A67:V67").Copy() | Out-Null
$ws.Range("A68:V68").PasteSpecial(-4122) | Out-Null

$ws.Range("A68").Value = 67
$ws.Range("B68").Value = "armenia"
$ws.Range("C68").Value = "premier-league"
$ws.Range("D68").Value = "2023-2024"
$ws.Range("E68").Value = 45232.45833333334
$ws.Range("F68").Value = "Pyunik Yerevan"
$ws.Range("G68").Value = 6
$ws.Range("H68").Value = "Van"
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1.11
$ws.Range("K68").Value = "31/10/2023 23:12"
$ws.Range("L68").Value = 1.11
$ws.Range("M68").Value = "02/11/2023 10:13"
$ws.Range("N68").Value = 8.6
$ws.Range("O68").Value = "31/10/2023 23:12"
$ws.Range("P68").Value = 9.949999999999999
$ws.Range("Q68").Value = "02/11/2023 10:56"
$ws.Range("R68").Value = 14.13
$ws.Range("S68").Value = "31/10/2023 23:12"
$ws.Range("T68").Value = 19.88
$ws.Range("U68").Value = "02/11/2023 10:56"
$ws.Range("V68").Value = "https://www.betexplorer.com/football/armenia/premier-league/pyunik-yerevan-van/GnXHmgCC/"

# --- Row 69 ---------------------------------------------------------------
$ws.Range("A67:V67").Copy() | Out-Null
$ws.Range("A69:V69").PasteSpecial(-4122) | Out-Null

$ws.Range("A69").Value = 68
$ws.Range("B69").Value = "armenia"
$ws.Range("C69").Value = "premier-league"
$ws.Range("D69").Value = "2023-2024"
$ws.Range("E69").Value = 45232.625
$ws.Range("F69").Value = "Ararat-Armenia"
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = "BKMA"
$ws.Range("I69").Value = 1
$ws.Range("J69").Value = 1.18
$ws.Range("K69").Value = "01/11/2023 03:12"
$ws.Range("L69").Value = 1.16
$ws.Range("M69").Value = "02/11/2023 14:54"
$ws.Range("N69").Value = 6.44
$ws.Range("O69").Value = "01/11/2023 03:12"
$ws.Range("P69").Value = 7.35
$ws.Range("Q69").Value = "02/11/2023 14:55"
$ws.Range("R69").Value = 10.88
$ws.Range("S69").Value = "01/11/2023 03:12"
$ws.Range("T69").Value = 16.18
$ws.Range("U69").Value = "02/11/2023 14:55"
$ws.Range("V69").Value = "https://www.betexplorer.com/football/armenia/premier-league/ararat-armenia-bkma/OCWLnDRI/"
